$p = $ppt.ActivePresentation

function Get-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            return $shp
        }
    }
    return $null
}

$newDateText = "12/8/2016"

# Slide master date placeholder
$m = $p.SlideMaster
$masterDateShape = Get-DatePlaceholder $m
if ($masterDateShape -ne $null) {
    $masterDateShape.TextFrame.TextRange.Text = $newDateText
}

# Every slide layout's date placeholder
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    $layoutDateShape = Get-DatePlaceholder $cl
    if ($layoutDateShape -ne $null) {
        $layoutDateShape.TextFrame.TextRange.Text = $newDateText
    }
}

# Slide 6 ("Käytetyt teknologiat"): add "Visual Studio" as the first bullet
# in the content placeholder, before "MarkedNet".
$s = $p.Slides.Item(6)
$contentShape = $s.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.InsertBefore("Visual Studio`r")
$newFirstPara = $tr.Paragraphs(1, 1)
$newFirstPara.LanguageID = "fi-FI"
